$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.250.49'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.225.95'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  -1.44%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '298.33'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -2.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '90.57'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -4.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.559'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -2.25%  '
$ws.Range('E9').Value = '  -5.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.20'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -5.06%  '
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.99'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -3.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.104'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('D14').Value = '2.565.56'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '2.230.57'
$ws.Range('E15').Value = '  -1.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.36'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.777'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -7.14%  '
$ws.Range('D18').Value = '44.091.73'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.21'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('E20').Value = '  -4.80%  '
$ws.Range('E21').Value = '  -5.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.31'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.54'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -0.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.82'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -4.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.84'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -6.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '39.11'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +2.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.21'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.37'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -4.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.26'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -3.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '151.47'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.50'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -8.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0762'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -4.41%  '
$ws.Range('E34').Value = '  -5.76%  '
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('E36').Value = '  -5.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.84'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -7.25%  '
$ws.Range('E38').Value = '  -8.50%  '
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.16'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -7.29%  '
$ws.Range('E41').Value = '  -4.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.44'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -10.08%  '
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').Value = '1.794.97'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.80'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +6.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.184'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -4.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '67.98'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '94.62'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -4.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '73.29'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -7.54%  '
$ws.Range('E50').Value = '  -4.68%  '
$ws.Range('E51').Value = '  -6.03%  '
